# Update cryptos list values per data refresh (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.811.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.097.81'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.83%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.96'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.389'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0780'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.50%  '

$ws.Range("E11").Value = '  +2.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.396.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.768'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.088.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.782.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("E24").Value = '  -1.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.22%  '

$ws.Range("E27").Value = '  +9.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.77%  '

$ws.Range("E29").Value = '  -1.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.14%  '

$ws.Range("E31").Value = '  +1.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0626'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.54%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.11%  '

$ws.Range("E37").Value = '  +4.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.66%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0997'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0214'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.456.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("E46").Value = '  +3.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.67%  '

$ws.Range("E50").Value = '  +2.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.291.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.80%  '
